$wb = $excel.ActiveWorkbook

# The report-for-handback update touches the two per-language sheets
# (zh-cn / de-de). For each one, the two tracked files (53b92d29-...md and
# ec025075-...md) have now been handed back: their status flips, the
# "Latest Target File" / "Latest Handback File" columns (E/F) get filled
# in (mirroring the existing "Latest Handoff" columns A/C, same display
# text + hyperlink target), and "Latest Handback DateTime" (G) moves from
# the 0001-01-01 placeholder to a real timestamp.

$sheets = @(
    @{
        Name = "zh-cn"
        HandbackTime = "2016-03-10 23:18:36"
        Row2 = @{
            TargetDisplay   = "53b92d29-16c3-441f-b557-35953c1e4d19.md"
            TargetUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/16e8fb4e99ce355c090d21ec564c8937f9d6d0d5/e2e/53b92d29-16c3-441f-b557-35953c1e4d19.md"
            HandbackDisplay = "53b92d29-16c3-441f-b557-35953c1e4d19.35c51cedde8ac13ae46980125ab00ec6f4ac640c.zh-cn.xlf"
            HandbackUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e28828a9c69a1ed50a56d3fa2be1d3bda7d5d48d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/53b92d29-16c3-441f-b557-35953c1e4d19.35c51cedde8ac13ae46980125ab00ec6f4ac640c.zh-cn.xlf"
        }
        Row3 = @{
            TargetDisplay   = "ec025075-cf11-4e78-a010-4af7558a3adc.md"
            TargetUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/16e8fb4e99ce355c090d21ec564c8937f9d6d0d5/e2e/ec025075-cf11-4e78-a010-4af7558a3adc.md"
            HandbackDisplay = "ec025075-cf11-4e78-a010-4af7558a3adc.8bfaa097c155c8b313cb91d081defebde8e1b3fb.zh-cn.xlf"
            HandbackUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e28828a9c69a1ed50a56d3fa2be1d3bda7d5d48d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ec025075-cf11-4e78-a010-4af7558a3adc.8bfaa097c155c8b313cb91d081defebde8e1b3fb.zh-cn.xlf"
        }
    },
    @{
        Name = "de-de"
        HandbackTime = "2016-03-10 23:18:54"
        Row2 = @{
            TargetDisplay   = "53b92d29-16c3-441f-b557-35953c1e4d19.md"
            TargetUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/16e8fb4e99ce355c090d21ec564c8937f9d6d0d5/e2e/53b92d29-16c3-441f-b557-35953c1e4d19.md"
            HandbackDisplay = "53b92d29-16c3-441f-b557-35953c1e4d19.35c51cedde8ac13ae46980125ab00ec6f4ac640c.de-de.xlf"
            HandbackUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/483e07b192e98807fbf5d1bd3b5792b4fc7706fb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/53b92d29-16c3-441f-b557-35953c1e4d19.35c51cedde8ac13ae46980125ab00ec6f4ac640c.de-de.xlf"
        }
        Row3 = @{
            TargetDisplay   = "ec025075-cf11-4e78-a010-4af7558a3adc.md"
            TargetUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/16e8fb4e99ce355c090d21ec564c8937f9d6d0d5/e2e/ec025075-cf11-4e78-a010-4af7558a3adc.md"
            HandbackDisplay = "ec025075-cf11-4e78-a010-4af7558a3adc.8bfaa097c155c8b313cb91d081defebde8e1b3fb.de-de.xlf"
            HandbackUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/483e07b192e98807fbf5d1bd3b5792b4fc7706fb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ec025075-cf11-4e78-a010-4af7558a3adc.8bfaa097c155c8b313cb91d081defebde8e1b3fb.de-de.xlf"
        }
    }
)

$handedBackStatus = "Handed back: in sync with en-US"
# Excel BGR-packed RGB, matching the workbook's existing hyperlink font color FF6495ED
$hyperlinkColor = 15570276

foreach ($sheetInfo in $sheets) {
    $ws = $wb.Worksheets.Item($sheetInfo.Name)

    # Status column updates for both tracked files
    $ws.Range("B2").Value = $handedBackStatus
    $ws.Range("B3").Value = $handedBackStatus

    # Row 2 (53b92d29-...md)
    $ws.Range("E2").Value = $sheetInfo.Row2.TargetDisplay
    $ws.Range("F2").Value = $sheetInfo.Row2.HandbackDisplay
    $ws.Range("G2").Value = $sheetInfo.HandbackTime

    $ws.Hyperlinks.Add($ws.Range("E2"), $sheetInfo.Row2.TargetUrl, "", "", $sheetInfo.Row2.TargetDisplay)
    $ws.Hyperlinks.Add($ws.Range("F2"), $sheetInfo.Row2.HandbackUrl, "", "", $sheetInfo.Row2.HandbackDisplay)

    # Row 3 (ec025075-...md)
    $ws.Range("E3").Value = $sheetInfo.Row3.TargetDisplay
    $ws.Range("F3").Value = $sheetInfo.Row3.HandbackDisplay
    $ws.Range("G3").Value = $sheetInfo.HandbackTime

    $ws.Hyperlinks.Add($ws.Range("E3"), $sheetInfo.Row3.TargetUrl, "", "", $sheetInfo.Row3.TargetDisplay)
    $ws.Hyperlinks.Add($ws.Range("F3"), $sheetInfo.Row3.HandbackUrl, "", "", $sheetInfo.Row3.HandbackDisplay)

    # Match the look of the existing hyperlink columns (A/C): underlined,
    # cornflower-blue text instead of the theme hyperlink style Add() uses.
    $ws.Range("E2:F3").Font.Underline = 2
    $ws.Range("E2:F3").Font.Color = $hyperlinkColor
}
